$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the TEXT_DATA table entry to BIOG_TEXT_DATA (new biography text table).
$found = $ws.Columns.Item(1).Find("TEXT_DATA")
$renameRow = $found.Row()
$ws.Cells.Item($renameRow, 1).Value = "BIOG_TEXT_DATA"

# 2. Append the newly-discovered tables at the bottom of the list (before sorting).
$lastRow = $ws.Cells.Item($ws.Rows.Count(), 1).End(-4162).Row()
$newTables = @("INDEXYEAR_TYPE_CODES", "FORMLABELS", "COPYMISSINGTABLES", "COPYTABLES", "COPYTABLESDEFAULT")
$r = $lastRow + 1
foreach ($name in $newTables) {
    $ws.Cells.Item($r, 1).Value = $name
    $ws.Cells.Item($r, 2).Formula = "=LOWER(A$r)"
    $ws.Cells.Item($r, 3).Formula = "=""ALTER TABLE ""&B$r&"" RENAME TO ""&A$r&"";"""
    $r = $r + 1
}
$lastRow = $r - 1

# 3. Re-sort the whole table (A2:C$lastRow) alphabetically by table name, same as the
#    original author did to keep the list easy to scan after adding new rows.
$sortRange = $ws.Range("A2:C" + $lastRow)
$sortKey = $ws.Range("A2:A" + $lastRow)
$sortRange.Sort($sortKey, 1)

# 4. Widen column A so the longer table names fit without wrapping as much.
$ws.Columns.Item(1).ColumnWidth = 38.33

# 5. Drop the now-unneeded extra row height on rows that used to wrap (column is wider now).
$ws.Rows.Item("1:" + $lastRow).EntireRow.AutoFit()

# 6. Reset the window scroll position / selection like the author left it.
$ws.Range("A67").Select()
